$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their exact text representation
# (e.g. "1.000", "0.4742") instead of being auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.458.61"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.54"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.97"
$ws.Range("E5").Value = "  +4.74%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4742"
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2910"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06509"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.94"
$ws.Range("E10").Value = "  +5.75%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07718"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.62"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7394"
$ws.Range("E13").Value = "  +9.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.880.83"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.125"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.32"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.471.61"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.42"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007541"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.128.07"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.238"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.184"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.266"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.44"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.83"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.935"
$ws.Range("E28").Value = "  +3.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1007"
$ws.Range("E29").Value = "  +2.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  +3.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.313"
$ws.Range("E32").Value = "  +3.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.103"
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04820"
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.124"
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6994"
$ws.Range("E36").Value = "  +2.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.724"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01861"
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.752"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.308"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.976"
$ws.Range("E42").Value = "  +5.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.94"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4201"
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8372"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.60"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.045"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.55"
$ws.Range("E50").Value = "  +4.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "915.01"
$ws.Range("E51").Value = "  -0.37%  "
